{"js": "const replacements = [\n  [\"2025-03-09 Sunday\", \"2025-03-10 Monday\"],\n  [\"647\u00d76=\", \"339\u00d73=\"],\n  [\"221\u00d76=\", \"546\u00d75=\"],\n  [\"562\u00d73=\", \"910\u00d73=\"],\n  [\"400\u00d72=\", \"584\u00d76=\"],\n  [\"260\u00d76=\", \"750\u00d79=\"],\n  [\"514\u00d73=\", \"413\u00d78=\"],\n  [\"927\u00d72=\", \"933\u00d79=\"],\n  [\"296\u00d74=\", \"789\u00d74=\"],\n  [\"728\u00d77=\", \"319\u00d72=\"],\n  [\"978\u00d72=\", \"518\u00d78=\"],\n  [\"173\u00d79=\", \"925\u00d77=\"],\n  [\"880\u00d75=\", \"787\u00d76=\"],\n  [\"435\u00d77=\", \"424\u00d78=\"],\n  [\"294\u00d72=\", \"422\u00d74=\"],\n  [\"493\u00d72=\", \"881\u00d78=\"],\n  [\"624\u00d73=\", \"805\u00d76=\"],\n  [\"530\u00d77=\", \"937\u00d79=\"],\n  [\"478\u00d74=\", \"516\u00d74=\"],\n  [\"638\u00d77=\", \"233\u00d78=\"],\n  [\"617\u00d74=\", \"674\u00d77=\"],\n  [\"649\u00d75=\", \"245\u00d76=\"],\n  [\"961\u00d76=\", \"366\u00d72=\"],\n  [\"374\u00d73=\", \"505\u00d78=\"],\n  [\"135\u00d75=\", \"312\u00d74=\"],\n  [\"863\u00d76=\", \"970\u00d77=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n  for (const range of results.items) {\n    range.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n  @(\"2025-03-09 Sunday\", \"2025-03-10 Monday\"),\n  @(\"647\u00d76=\", \"339\u00d73=\"),\n  @(\"221\u00d76=\", \"546\u00d75=\"),\n  @(\"562\u00d73=\", \"910\u00d73=\"),\n  @(\"400\u00d72=\", \"584\u00d76=\"),\n  @(\"260\u00d76=\", \"750\u00d79=\"),\n  @(\"514\u00d73=\", \"413\u00d78=\"),\n  @(\"927\u00d72=\", \"933\u00d79=\"),\n  @(\"296\u00d74=\", \"789\u00d74=\"),\n  @(\"728\u00d77=\", \"319\u00d72=\"),\n  @(\"978\u00d72=\", \"518\u00d78=\"),\n  @(\"173\u00d79=\", \"925\u00d77=\"),\n  @(\"880\u00d75=\", \"787\u00d76=\"),\n  @(\"435\u00d77=\", \"424\u00d78=\"),\n  @(\"294\u00d72=\", \"422\u00d74=\"),\n  @(\"493\u00d72=\", \"881\u00d78=\"),\n  @(\"624\u00d73=\", \"805\u00d76=\"),\n  @(\"530\u00d77=\", \"937\u00d79=\"),\n  @(\"478\u00d74=\", \"516\u00d74=\"),\n  @(\"638\u00d77=\", \"233\u00d78=\"),\n  @(\"617\u00d74=\", \"674\u00d77=\"),\n  @(\"649\u00d75=\", \"245\u00d76=\"),\n  @(\"961\u00d76=\", \"366\u00d72=\"),\n  @(\"374\u00d73=\", \"505\u00d78=\"),\n  @(\"135\u00d75=\", \"312\u00d74=\"),\n  @(\"863\u00d76=\", \"970\u00d77=\")\n)\n\nforeach ($pair in $replacements) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.ClearFormatting()\n  $find.Replacement.Text = $newText\n  $find.Execute($oldText, $false, $false, $false, $false, $false, $true, 0, $false, $newText, 2)\n}\n"}
